$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blackbox")
$rng = $ws.Range("G3:G8")
$rng.Style = "Good"
